$wb = $excel.ActiveWorkbook

# Sheets in this workbook: 1=model, 2=settings, 3=survey, 4=properties
$propsSheet = $wb.Worksheets.Item("properties")

# Add two new rows of data to the "properties" sheet
$propsSheet.Range("A5").Value = "FormType"
$propsSheet.Range("B5").Value = "default"
$propsSheet.Range("C5").Value = "FormType.formType"
$propsSheet.Range("D5").Value = "string"
$propsSheet.Range("E5").Value = "SURVEY"

$propsSheet.Range("A6").Value = "SurveyUtil"
$propsSheet.Range("B6").Value = "default"
$propsSheet.Range("C6").Value = "SurveyUtil.formId"
$propsSheet.Range("D6").Value = "string"
$propsSheet.Range("E6").Value = "wrong_form"

# Give the new rows the same (plain, non-centered) formatting used by columns A-D of the
# existing rows -- applied across the whole A:E range for rows 5-6.
$propsSheet.Range("A2").Copy() | Out-Null
$propsSheet.Range("A5:E6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Switch the selected/active tab from "model" to "properties", with a new cell selection.
$propsSheet.Activate()
$propsSheet.Range("F11").Select() | Out-Null
